# Weekly update: insert a new weekly price record as row 31 on the
# "Hortaliza - Acelga" sheet, shifting the existing rows (old 31..62)
# down to (32..63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31; everything below (including
# the old row 31) shifts down by one.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44671
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112009
$ws.Range("G31").Value = "Acelga"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 1300
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1400
$ws.Range("N31").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 467
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = "Hortaliza"
